# Append additional process-log rows (4-15) to the sheet, mirroring the
# existing rows 2-3, and move the trailing empty "Error Details" marker
# cell down to the new last row (H15).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the stray empty marker left on the previous last row (H3) - the
# row no longer ends the table, so it no longer carries the marker cell.
$ws.Range("H3").Value = ""

$rows = @(
    @("5/30/2022 15:44", "Monday", "Alcon", "Processed", "Sent", "Alcon_Status Report_05302022.xlsx", "Alcon_Document Expiration Report_05302022.xlsx"),
    @("5/30/2022 15:52", "Monday", "Alcon", "Processed", "Sent", "Alcon_Status Report_05302022.xlsx", "Alcon_Document Expiration Report_05302022.xlsx"),
    @("5/30/2022 16:26", "Monday", "Alcon", "Processed", "Sent", "Alcon_Status Report_05302022.xlsx", "Alcon_Document Expiration Report_05302022.xlsx"),
    @("5/30/2022 16:34", "Monday", "Alcon", "Processed", "Sent", "Alcon_Status Report_05302022.xlsx", "Alcon_Document Expiration Report_05302022.xlsx"),
    @("5/30/2022 16:40", "Monday", "Alcon", "Processed", "Sent", "Alcon_Status Report_05302022.xlsx", "Alcon_Document Expiration Report_05302022.xlsx"),
    @("5/30/2022 16:58", "Monday", "Alcon", "Processed", "Sent", "Alcon_Status Report_05302022.xlsx", "Alcon_Document Expiration Report_05302022.xlsx"),
    @("5/30/2022 17:04", "Monday", "Alcon", "Processed", "Sent", "Alcon_Status Report_05302022.xlsx", "Alcon_Document Expiration Report_05302022.xlsx"),
    @("5/30/2022 17:10", "Monday", "Alcon", "Processed", "Sent", "Alcon_Status Report_05302022.xlsx", "Alcon_Document Expiration Report_05302022.xlsx"),
    @("5/30/2022 19:33", "Monday", "Alcon", "Processed", "Sent", "Alcon_Status Report_05302022.xlsx", "Alcon_Document Expiration Report_05302022.xlsx"),
    @("5/30/2022 19:49", "Monday", "Alcon", "Processed", "Sent", "Alcon_Status Report_05302022.xlsx", "Alcon_Document Expiration Report_05302022.xlsx"),
    @("5/30/2022 19:55", "Monday", "Alcon", "Processed", "Sent", "Alcon_Status Report_05302022.xlsx", "Alcon_Document Expiration Report_05302022.xlsx"),
    @("5/30/2022 20:02", "Monday", "Alcon", "Processed", "Sent", "Alcon_Status Report_05302022.xlsx", "Alcon_Document Expiration Report_05302022.xlsx")
)

$startRow = 4
for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).Value = $data[3]
    $ws.Cells.Item($r, 5).Value = $data[4]
    $ws.Cells.Item($r, 6).Value = $data[5]
    $ws.Cells.Item($r, 7).Value = $data[6]
}

# New last row (15) carries the empty "Error Details" marker cell, matching
# the pattern previously seen on the old last row (3). A plain "" assignment
# clears the cell outright instead of leaving a present-but-empty value, so
# use the leading-apostrophe text-prefix trick to force a stored empty
# string, then reset the style so it doesn't pick up the quote-prefix flag.
$ws.Cells.Item(15, 8).Value = "'"
$ws.Cells.Item(15, 8).Style = "Normal"
